$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 217, pushing the existing rows 217-223 down to 218-224.
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the weekly entry.
$ws.Range("A217").Value = 9
$ws.Range("B217").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C217").Value = "Metropolitana"
$ws.Range("D217").Value = 44509
$ws.Range("E217").Value = 13
$ws.Range("F217").Value = 100112052
$ws.Range("G217").Value = "Albahaca"
$ws.Range("H217").Value = "Sin especificar"
$ws.Range("I217").Value = "Primera"
$ws.Range("J217").Value = 79
$ws.Range("K217").Value = 5000
$ws.Range("L217").Value = 6000
$ws.Range("M217").Value = 5506
$ws.Range("N217").Value = "`$/docena de matas"
$ws.Range("O217").Value = "Provincia de Chacabuco"
$ws.Range("P217").Value = 918
$ws.Range("Q217").Value = 6
$ws.Range("R217").Value = "Hortaliza"
